# Re-run of the NATMI cell-cell-signalling export with updated TPM input
# (commit: "update scripts wuth new tpm"). The "Sending cluster" vs
# "Target cluster" expression numbers used to build the edge-weight /
# specificity columns (E:T) for each Ligand-Receptor row were
# recalculated from the refreshed TPM matrix. Sending/Receptor cluster
# labels (A:D) are unchanged; only the numeric columns below them move.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.344207
$ws.Range("H2").Value = 1.032621
$ws.Range("I2").Value = 0.6985282229833164
$ws.Range("J2").Value = 0.6985282229833165
$ws.Range("M2").Value = 0.9519303333333333
$ws.Range("N2").Value = 2.855791
$ws.Range("O2").Value = 0.2529401125772161
$ws.Range("P2").Value = 0.2529401125772162
$ws.Range("Q2").Value = 0.3276610842456666
$ws.Range("R2").Value = 2.948949758211
$ws.Range("S2").Value = 0.1766858073597628
$ws.Range("T2").Value = 0.1766858073597628

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.344207
$ws.Range("H3").Value = 1.032621
$ws.Range("I3").Value = 0.6985282229833164
$ws.Range("J3").Value = 0.6985282229833165
$ws.Range("O3").Value = 0.4278029870374648
$ws.Range("P3").Value = 0.4278029870374649
$ws.Range("Q3").Value = 0.5541801541400001
$ws.Range("R3").Value = 4.987621387260001
$ws.Range("S3").Value = 0.298832460322235
$ws.Range("T3").Value = 0.2988324603222351

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.344207
$ws.Range("H4").Value = 1.032621
$ws.Range("I4").Value = 0.6985282229833164
$ws.Range("J4").Value = 0.6985282229833165
$ws.Range("M4").Value = 0.9848326666666667
$ws.Range("N4").Value = 2.954498
$ws.Range("O4").Value = 0.261682685017622
$ws.Range("P4").Value = 0.2616826850176221
$ws.Range("Q4").Value = 0.3389862976953333
$ws.Range("R4").Value = 3.050876679258
$ws.Range("S4").Value = 0.1827927409508625
$ws.Range("T4").Value = 0.1827927409508625

# Row 5: ECs -> Resolving-Mac
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.344207
$ws.Range("H5").Value = 1.032621
$ws.Range("I5").Value = 0.6985282229833164
$ws.Range("J5").Value = 0.6985282229833165
$ws.Range("M5").Value = 0.2166783333333333
$ws.Range("N5").Value = 0.6500349999999999
$ws.Range("O5").Value = 0.05757421536769695
$ws.Range("P5").Value = 0.05757421536769697
$ws.Range("Q5").Value = 0.07458219908166665
$ws.Range("R5").Value = 0.6712397917349999
$ws.Range("S5").Value = 0.0402172143504561
$ws.Range("T5").Value = 0.04021721435045612

# Row 6: FAPs -> ECs
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1485533333333333
$ws.Range("H6").Value = 0.44566
$ws.Range("I6").Value = 0.3014717770166836
$ws.Range("J6").Value = 0.3014717770166836
$ws.Range("M6").Value = 0.9519303333333333
$ws.Range("N6").Value = 2.855791
$ws.Range("O6").Value = 0.2529401125772161
$ws.Range("P6").Value = 0.2529401125772162
$ws.Range("Q6").Value = 0.1414124241177778
$ws.Range("R6").Value = 1.27271181706
$ws.Range("S6").Value = 0.07625430521745334
$ws.Range("T6").Value = 0.07625430521745336

# Row 7: FAPs -> FAPs
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1485533333333333
$ws.Range("H7").Value = 0.44566
$ws.Range("I7").Value = 0.3014717770166836
$ws.Range("J7").Value = 0.3014717770166836
$ws.Range("O7").Value = 0.4278029870374648
$ws.Range("P7").Value = 0.4278029870374649
$ws.Range("Q7").Value = 0.2391738377333334
$ws.Range("R7").Value = 2.1525645396
$ws.Range("S7").Value = 0.1289705267152298
$ws.Range("T7").Value = 0.1289705267152298

# Row 8: FAPs -> MuSCs
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1485533333333333
$ws.Range("H8").Value = 0.44566
$ws.Range("I8").Value = 0.3014717770166836
$ws.Range("J8").Value = 0.3014717770166836
$ws.Range("M8").Value = 0.9848326666666667
$ws.Range("N8").Value = 2.954498
$ws.Range("O8").Value = 0.261682685017622
$ws.Range("P8").Value = 0.2616826850176221
$ws.Range("Q8").Value = 0.1463001754088889
$ws.Range("R8").Value = 1.31670157868
$ws.Range("S8").Value = 0.0788899440667596
$ws.Range("T8").Value = 0.07888994406675962

# Row 9: FAPs -> Resolving-Mac
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1485533333333333
$ws.Range("H9").Value = 0.44566
$ws.Range("I9").Value = 0.3014717770166836
$ws.Range("J9").Value = 0.3014717770166836
$ws.Range("M9").Value = 0.2166783333333333
$ws.Range("N9").Value = 0.6500349999999999
$ws.Range("O9").Value = 0.05757421536769695
$ws.Range("P9").Value = 0.05757421536769697
$ws.Range("Q9").Value = 0.03218828867777777
$ws.Range("R9").Value = 0.2896945981
$ws.Range("S9").Value = 0.01735700101724085
$ws.Range("T9").Value = 0.01735700101724086
